# Updated cryptos list on Fri Sep 22 23:23:44 UTC 2023 with GitHub Actions
#
# This script refreshes the crypto price table on Sheet1. Most rows keep
# the same coin in place and simply get a refreshed Price (column D) and/or
# Volume(1h) (column E) figure. A handful of rows were re-ranked by the
# source feed, so those rows' Coin (B), Link (C), Price (D) and Volume (E)
# are rewritten wholesale (this covers the WrappedEther/Polkadot swap at
# rows 13-14, the Aave/MXToken swap at rows 42-43, and the insertion of
# BabyDogeCoin ahead of Algorand/Cronos - which also pushes EnergySwap off
# the bottom of the 50-row list).
#
# NOTE: column D ("Price") stores plain text, not numbers (e.g. "211.40",
# "1.00", "26.626.78"), and some values would otherwise be auto-detected
# and coerced into floating point numbers by Excel (losing trailing zeros
# or precision, e.g. "211.40" -> 211.4). To keep them as text exactly as
# written, values that parse as a plain number are assigned with a leading
# apostrophe (forces text entry, like typing into Excel by hand) and then
# the cell style is reset back to "Normal" so no stray numeric/quote-prefix
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin: price + volume refresh (price has two dots, stays text naturally)
$ws.Cells.Item(2, 4).Value = "26.626.78"
$ws.Cells.Item(2, 5).Value = "  -0.02%  "

# Row 3 - Ethereum: price refresh only (price has two dots, stays text naturally)
$ws.Cells.Item(3, 4).Value = "1.596.69"

# Row 4 - TetherUSD: volume refresh only
$ws.Cells.Item(4, 5).Value = "  +0.13%  "

# Row 5 - BNB: price + volume refresh
$ws.Cells.Item(5, 4).Value = "'211.40"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.17%  "

# Row 6 - XRP: price + volume refresh
$ws.Cells.Item(6, 4).Value = "'0.515"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.13%  "

# Row 7 - USDC: volume refresh only
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

# Row 8 - Dogecoin: volume refresh only
$ws.Cells.Item(8, 5).Value = "  +0.15%  "

# Row 9 - Cardano: volume refresh only
$ws.Cells.Item(9, 5).Value = "  -0.29%  "

# Row 10 - Solana: price + volume refresh
$ws.Cells.Item(10, 4).Value = "'19.52"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.30%  "

# Row 11 - TRON: volume refresh only
$ws.Cells.Item(11, 5).Value = "  +0.48%  "

# Row 12 - WrappedliquidstakedEther2.0: price + volume refresh (two dots, stays text naturally)
$ws.Cells.Item(12, 4).Value = "1.821.19"
$ws.Cells.Item(12, 5).Value = "  +0.59%  "

# Row 13: was WrappedEther, now Polkadot (re-ranked)
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "'4.03"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.08%  "

# Row 14: was Polkadot, now WrappedEther (re-ranked) (two dots, stays text naturally)
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.552.10"
$ws.Cells.Item(14, 5).Value = "  -2.13%  "

# Row 15 - Polygon: volume refresh only
$ws.Cells.Item(15, 5).Value = "  -0.04%  "

# Row 16 - Litecoin: price + volume refresh
$ws.Cells.Item(16, 4).Value = "'64.37"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.58%  "

# Row 17 - WrappedBTC: price refresh only (two dots, stays text naturally)
$ws.Cells.Item(17, 4).Value = "26.617.51"

# Row 18 - ShibaInu: volume refresh only
$ws.Cells.Item(18, 5).Value = "  +0.51%  "

# Row 19 - BitcoinCash: price + volume refresh
$ws.Cells.Item(19, 4).Value = "'208.50"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.37%  "

# Row 20 - Dai: volume refresh only
$ws.Cells.Item(20, 5).Value = "  +0.07%  "

# Row 21 - Chainlink: price + volume refresh
$ws.Cells.Item(21, 4).Value = "'6.95"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.55%  "

# Row 22 - Uniswap: volume refresh only
$ws.Cells.Item(22, 5).Value = "  +0.28%  "

# Row 23 - Toncoin: volume refresh only
$ws.Cells.Item(23, 5).Value = "  -2.34%  "

# Row 24 - Avalanche: volume refresh only
$ws.Cells.Item(24, 5).Value = "  +0.27%  "

# Row 25 - Monero: price + volume refresh
$ws.Cells.Item(25, 4).Value = "'145.01"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.71%  "

# Row 26 - BinanceUSD: volume refresh only
$ws.Cells.Item(26, 5).Value = "  +0.08%  "

# Row 29 - EthereumClassic: price + volume refresh
$ws.Cells.Item(29, 4).Value = "'15.24"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.30%  "

# Row 30 - Hedera: price + volume refresh
$ws.Cells.Item(30, 4).Value = "'0.0506"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.16%  "

# Row 31 - PancakeSwap: volume refresh only
$ws.Cells.Item(31, 5).Value = "  +0.36%  "

# Row 32 - Filecoin: volume refresh only
$ws.Cells.Item(32, 5).Value = "  -0.03%  "

# Row 33 - ImmutableX: volume refresh only
$ws.Cells.Item(33, 5).Value = "  -3.50%  "

# Row 34 - InternetComputer(DFINITY): price + volume refresh
$ws.Cells.Item(34, 4).Value = "'2.92"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.58%  "

# Row 35 - Maker: price + volume refresh (two dots, stays text naturally)
$ws.Cells.Item(35, 4).Value = "1.275.93"
$ws.Cells.Item(35, 5).Value = "  -2.71%  "

# Row 36 - HuobiToken: volume refresh only
$ws.Cells.Item(36, 5).Value = "  +0.37%  "

# Row 37 - LidoDAOToken: volume refresh only
$ws.Cells.Item(37, 5).Value = "  +0.79%  "

# Row 38 - VeChain: volume refresh only
$ws.Cells.Item(38, 5).Value = "  -0.61%  "

# Row 39 - ARBITRUM: price + volume refresh
$ws.Cells.Item(39, 4).Value = "'0.843"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.99%  "

# Row 40 - PaxDollar: volume refresh only
$ws.Cells.Item(40, 5).Value = "  +0.08%  "

# Row 41 - FraxShare: price + volume refresh
$ws.Cells.Item(41, 4).Value = "'5.47"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.66%  "

# Row 42: was Aave, now MXToken (re-ranked)
$ws.Cells.Item(42, 2).Value = "MXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(42, 4).Value = "'2.20"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.90%  "

# Row 43: was MXToken, now Aave (re-ranked)
$ws.Cells.Item(43, 2).Value = "Aave"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(43, 4).Value = "'64.44"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.83%  "

# Row 44 - TrustWalletToken: volume refresh only
$ws.Cells.Item(44, 5).Value = "  -0.63%  "

# Row 45 - WEMIXToken: price + volume refresh
$ws.Cells.Item(45, 4).Value = "'0.920"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +9.39%  "

# Row 46 - RocketPoolETH: price + volume refresh (two dots, stays text naturally)
$ws.Cells.Item(46, 4).Value = "1.733.41"
$ws.Cells.Item(46, 5).Value = "  +0.54%  "

# Row 47 - Quant: price + volume refresh
$ws.Cells.Item(47, 4).Value = "'89.89"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.61%  "

# Row 48 - RenderToken: volume refresh only
$ws.Cells.Item(48, 5).Value = "  -0.26%  "

# Row 49: BabyDogeCoin newly ranked ahead of Algorand (was Algorand).
# Price uses a Unicode subscript-six (U+2086) digit-grouping character,
# which must be built via [char] + string-subexpression interpolation
# (direct string concatenation with + on a [char] coerces to numeric
# addition in this engine, and `u{...} escapes are not supported). The
# resulting text ("0.0<sub6>0104") does not parse as a number, so no
# apostrophe/style trick is required here.
$subscriptSix = [char]0x2086
$babyDogePrice = "0.0$($subscriptSix)0104"
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = $babyDogePrice
$ws.Cells.Item(49, 5).Value = "  -1.54%  "

# Row 50: was Cronos, now Algorand (shifted down by the BabyDogeCoin insert)
$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).Value = "'0.102"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +4.31%  "

# Row 51: was EnergySwap, now Cronos (EnergySwap drops off the bottom of the list)
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.0507"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.43%  "
